$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove C2 entirely (was an erroneous naive-forecaster value)
$ws.Range("C2").ClearContents()

# Update recalculated forecast values (precision/bugfix corrections)
$ws.Range("C5").Value = 1.891565607550105
$ws.Range("E5").Value = 1.609625625600009

$ws.Range("C6").Value = 1.113165545862116
$ws.Range("E6").Value = 1.609625625599986

$ws.Range("E7").Value = 3.238605209600021

$ws.Range("C8").Value = 1.384186838979828
$ws.Range("E8").Value = 2.777885851461526

$ws.Range("C9").Value = 2.349355943833076

$ws.Range("C10").Value = 1.78642563555842

$ws.Range("E11").Value = 1.216098605743343

$ws.Range("C12").Value = 1.282262557986469
$ws.Range("E12").Value = 1.784618024189033

$ws.Range("C13").Value = 2.247109253368307

$ws.Range("C14").Value = -4.247034401476779
$ws.Range("E14").Value = -12.19860234240002

$ws.Range("E16").Value = -0.5376914776811237

$ws.Range("E17").Value = -4.829433539906869

$ws.Range("C18").Value = -0.244366674180263
$ws.Range("E18").Value = -1.64927836088965
